# Auto-generated Excel COM-interop script
# Applies cached market-price / profit recalculation updates to the
# Chocobo_Profits workbook (per scheduled-runner data refresh).
#
# For each touched Leve row we overwrite the cached market-price /
# profit columns (H:N) with freshly pulled values. A few rows also
# gain or lose a cell entirely (e.g. LeveProfitNQ/HQ becoming
# inapplicable once HQ pricing data disappears, or appearing once it
# becomes available) - those use ClearContents()/Value as appropriate.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 318.2
$ws.Range("I33").Value = 264.66666
$ws.Range("K33").Value = 264.66666
$ws.Range("M33").Value = -35.66665999999998
$ws.Range("H38").Value = 6125.143
$ws.Range("I38").Value = 575.2
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 1725.6
$ws.Range("L38").Value = 60000
$ws.Range("M38").Value = -1353.6
$ws.Range("N38").Value = -60744
$ws.Range("H51").Value = 50007500
$ws.Range("I51").Value = 100000000
$ws.Range("J51").Value = 15000
$ws.Range("K51").Value = 100000000
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = -99999516
$ws.Range("N51").Value = -15968
$ws.Range("H111").Value = 2766
$ws.Range("J111").Value = 2766
$ws.Range("L111").Value = 8298
$ws.Range("N111").Value = -14432
$ws.Range("H123").Value = 41835
$ws.Range("J123").Value = 41835
$ws.Range("L123").Value = 41835
$ws.Range("N123").Value = -51635
$ws.Range("H135").Value = 894
$ws.Range("I135").Value = 862.7143
$ws.Range("J135").Value = 967
$ws.Range("K135").Value = 7764.428699999999
$ws.Range("L135").Value = 8703
$ws.Range("M135").Value = -5229.428699999999
$ws.Range("N135").Value = -13773
$ws.Range("H137").Value = 853845.75
$ws.Range("I137").Value = 2271681.5
$ws.Range("J137").Value = 3144.2856
$ws.Range("K137").Value = 6815044.5
$ws.Range("L137").Value = 9432.856800000001
$ws.Range("M137").Value = -6812494.5
$ws.Range("N137").Value = -14532.8568
$ws.Range("H138").Value = 3670.7144
$ws.Range("I138").Value = 3127.2856
$ws.Range("J138").Value = 4214.143
$ws.Range("K138").Value = 9381.856800000001
$ws.Range("L138").Value = 12642.429
$ws.Range("M138").Value = -4241.856800000001
$ws.Range("N138").Value = -22922.429

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3487.7415
$ws.Range("I32").Value = 3369.36
$ws.Range("J32").Value = 4121.9287
$ws.Range("K32").Value = 3369.36
$ws.Range("L32").Value = 4121.9287
$ws.Range("M32").Value = -3082.36
$ws.Range("N32").Value = -4695.9287
$ws.Range("H35").Value = 18099.8
$ws.Range("J35").Value = 35999.5
$ws.Range("L35").Value = 35999.5
$ws.Range("N35").Value = -36811.5
$ws.Range("H61").Value = 3602.4
$ws.Range("I61").Value = 3602.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3602.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3390.4
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 267761.12
$ws.Range("I74").Value = 561642.2
$ws.Range("J74").Value = 1868.762
$ws.Range("K74").Value = 561642.2
$ws.Range("L74").Value = 1868.762
$ws.Range("M74").Value = -560768.2
$ws.Range("N74").Value = -3616.762
$ws.Range("H77").Value = 267761.12
$ws.Range("I77").Value = 561642.2
$ws.Range("J77").Value = 1868.762
$ws.Range("K77").Value = 2808211
$ws.Range("L77").Value = 9343.809999999999
$ws.Range("M77").Value = -2803843
$ws.Range("N77").Value = -18079.81
$ws.Range("H122").Value = 3270.5386
$ws.Range("I122").Value = 3057.7856
$ws.Range("J122").Value = 3812.0908
$ws.Range("K122").Value = 9173.356800000001
$ws.Range("L122").Value = 11436.2724
$ws.Range("M122").Value = -6723.356800000001
$ws.Range("N122").Value = -16336.2724
$ws.Range("H132").Value = 2653.9534
$ws.Range("I132").Value = 2274.7812
$ws.Range("J132").Value = 3757
$ws.Range("K132").Value = 6824.3436
$ws.Range("L132").Value = 11271
$ws.Range("M132").Value = -4294.3436
$ws.Range("N132").Value = -16331
$ws.Range("H136").Value = 3602.4
$ws.Range("I136").Value = 3602.4
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10807.2
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8257.200000000001
$ws.Range("N136").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 284353.6
$ws.Range("I31").Value = 1502910.8
$ws.Range("J31").Value = 3148.1025
$ws.Range("K31").Value = 1502910.8
$ws.Range("L31").Value = 3148.1025
$ws.Range("M31").Value = -1502615.8
$ws.Range("N31").Value = -3738.1025
$ws.Range("H34").Value = 284353.6
$ws.Range("I34").Value = 1502910.8
$ws.Range("J34").Value = 3148.1025
$ws.Range("K34").Value = 1502910.8
$ws.Range("L34").Value = 3148.1025
$ws.Range("M34").Value = -1502708.8
$ws.Range("N34").Value = -3552.1025
$ws.Range("H58").Value = 2500.611
$ws.Range("I58").Value = 1442.08
$ws.Range("J58").Value = 4906.364
$ws.Range("K58").Value = 1442.08
$ws.Range("L58").Value = 4906.364
$ws.Range("M58").Value = -1239.08
$ws.Range("N58").Value = -5312.364
$ws.Range("H132").Value = 3499.8823
$ws.Range("I132").Value = 2088.4
$ws.Range("J132").Value = 5516.2856
$ws.Range("K132").Value = 6265.200000000001
$ws.Range("L132").Value = 16548.8568
$ws.Range("M132").Value = -3735.200000000001
$ws.Range("N132").Value = -21608.8568
$ws.Range("H134").Value = 8174.1875
$ws.Range("I134").Value = 9393.166999999999
$ws.Range("J134").Value = 4517.25
$ws.Range("K134").Value = 28179.501
$ws.Range("L134").Value = 13551.75
$ws.Range("M134").Value = -25644.501
$ws.Range("N134").Value = -18621.75
$ws.Range("H136").Value = 2500.611
$ws.Range("I136").Value = 1442.08
$ws.Range("J136").Value = 4906.364
$ws.Range("K136").Value = 4326.24
$ws.Range("L136").Value = 14719.092
$ws.Range("M136").Value = -1776.24
$ws.Range("N136").Value = -19819.092

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 3170
$ws.Range("J25").Value = 3170
$ws.Range("L25").Value = 9510
$ws.Range("N25").Value = -9848
$ws.Range("H29").Value = 299.66666
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 299.66666
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 898.9999799999999
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -1452.99998
$ws.Range("H30").Value = 3170
$ws.Range("J30").Value = 3170
$ws.Range("L30").Value = 9510
$ws.Range("N30").Value = -9714
$ws.Range("H36").Value = 2000
$ws.Range("J36").Value = 2000
$ws.Range("L36").Value = 6000
$ws.Range("N36").Value = -6338
$ws.Range("H37").Value = 58895350
$ws.Range("J37").Value = 58895350
$ws.Range("L37").Value = 176686050
$ws.Range("N37").Value = -176686274
$ws.Range("H124").Value = 7176.6665
$ws.Range("I124").Value = 5765
$ws.Range("J124").Value = 10000
$ws.Range("K124").Value = 17295
$ws.Range("L124").Value = 30000
$ws.Range("M124").Value = -12385
$ws.Range("N124").Value = -39820

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4233.25
$ws.Range("I126").Value = 3411.1272
$ws.Range("J126").Value = 5603.4546
$ws.Range("K126").Value = 10233.3816
$ws.Range("L126").Value = 16810.3638
$ws.Range("M126").Value = -7763.381600000001
$ws.Range("N126").Value = -21750.3638
$ws.Range("H132").Value = 2751
$ws.Range("I132").Value = 2293
$ws.Range("J132").Value = 3636.4666
$ws.Range("K132").Value = 6879
$ws.Range("L132").Value = 10909.3998
$ws.Range("M132").Value = -4349
$ws.Range("N132").Value = -15969.3998

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 130.6
$ws.Range("I55").Value = 100.333336
$ws.Range("J55").Value = 176
$ws.Range("K55").Value = 100.333336
$ws.Range("L55").Value = 176
$ws.Range("M55").Value = 72.666664
$ws.Range("N55").Value = -522

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 30300
$ws.Range("J108").Value = 30300
$ws.Range("L108").Value = 30300
$ws.Range("N108").Value = -37980
$ws.Range("H136").Value = 3159.2163
$ws.Range("I136").Value = 1468.4445
$ws.Range("J136").Value = 4761
$ws.Range("K136").Value = 4405.333500000001
$ws.Range("L136").Value = 14283
$ws.Range("M136").Value = -1855.333500000001
$ws.Range("N136").Value = -19383

